# Update the workbook to reflect the latest data refresh:
#  - The "Förändrad" (changed) date in column C moves from 2023-10-06 (45205)
#    to 2023-10-07 (45206) for every existing data row (rows 2-251).
#  - Row 251 gains an explicit row height (matches the new row below it).
#  - A new data row (252) is appended for case "A 48177-2023".
#  - The used-range dimension grows from A1:Y251 to A1:Y252 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 251

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45206
}

# Row 251 now has an explicit custom row height (15pt), like the new last row.
$ws.Rows.Item($lastRow).RowHeight = 15

# Append the new record as row 252.
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "A 48177-2023"

$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow, 2).Value = 45205

$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow, 3).Value = 45206

$ws.Cells.Item($newRow, 4).Value = "VÄSTMANLANDS LÄN"
$ws.Cells.Item($newRow, 5).Value = "NORBERG"
$ws.Cells.Item($newRow, 6).Value = "Övriga Aktiebolag"

$ws.Cells.Item($newRow, 7).Value = 18.4

for ($c = 8; $c -le 17; $c++) {
    $ws.Cells.Item($newRow, $c).Value = 0
}

# Column R keeps the same wrapped, empty-text style used throughout the sheet.
$ws.Cells.Item($newRow, 18).WrapText = $true
$ws.Cells.Item($newRow, 18).Value = ""

Write-Host "Updated $($lastRow - 1) rows and appended row $newRow"
